$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "FILE_NAME"
$ws.Range("B1").Value = "FILE_TYPE"
$ws.Range("C1").Value = "BUCKET_LOCATION"
$ws.Range("D1").Value = "FILE_TAG_1"
$ws.Range("E1").Value = "FILE_TAG_2"
$ws.Range("F1").Value = "FILE_TAG_3"
$ws.Range("G1").Value = "FILE_TAG_4"
$ws.Range("H1").Value = "FILE_TAG_5"

# Row 2
$ws.Range("A2").Value = "AWS-Achieves_FED-Ramp-JPEG-2.jpg"
$ws.Range("B2").Value = "jpg"
$ws.Range("C2").Value = "uconn-sdp-team11-unprocessed-docs"
$ws.Range("D2").Value = "the"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""

# Row 3
$ws.Range("A3").Value = "AWS-Achieves_FED-Ramp-JPEG.jpg"
$ws.Range("B3").Value = "jpg"
$ws.Range("C3").Value = "uconn-sdp-team11-unprocessed-docs"
$ws.Range("D3").Value = "a"
$ws.Range("E3").Value = "for"
$ws.Range("F3").Value = "is"
$ws.Range("G3").Value = "the"
$ws.Range("H3").Value = ""

# Row 4 - remove entirely
$ws.Range("A4:H4").Value = ""
$ws.Rows.Item(4).Delete()
